# Switch the deck's active design from "Integral" (Red Violet) back to the
# built-in "Office Theme" colour scheme.
#
# The presentation's single slide master (and therefore every slide, layout
# and the presentation's default theme) is backed by one DrawingML theme
# part. Its 12-colour scheme currently holds the "Red Violet" palette used
# by the "Integral" design; we repoint every slot to the standard "Office"
# palette, which is what the author's edit applied.
#
# PowerPoint exposes the live colour scheme through Slide.ThemeColorScheme
# (shared by every slide because they all inherit the one slide master), so
# editing it from slide 1 updates the design used throughout the deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink),
# expressed as OLE/VBA RGB() integers (0xBBGGRR).
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
